$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (Jan-16 run -> shifts to D, Issues -> shifts to E)
$ws.Columns.Item(3).Insert()

# Re-order / extend the bank rows: HSBC then NDB, then the three new banks,
# each with a sequential "#".
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "HSBC"
$ws.Range("D7").Value = "Didn’t take merchant name of first two promos in each category"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "NDB"
$ws.Range("D8").Value = "Ok"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "People Bank"
$ws.Range("D9").Value = "Ok"

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "Sampath"
$ws.Range("D10").Value = "Ok"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Seylan"
$ws.Range("D11").Value = "Ok"
$ws.Range("E11").Value = "Need to give proper summary in console"

# Fill in the "Ok" result for the Feb 4th run on the already-present banks
$ws.Range("C2").Value = "Ok"

# New column header for the latest run (typed last, so it lands at the end
# of the shared-string table)
$ws.Range("C1").Value = "Feb 4th run"

# Bold the header row
$ws.Rows.Item(1).Font.Bold = $true

# Selection / page setup to mirror the saved state
$ws.Range("E12").Select()
$ws.PageSetup.Orientation = 1
